# Update "想去人数" (want-to-go count) figures that changed between data
# refreshes for the gh-pages generated output.
#
# Sheet "展览" (exhibitions) and sheet "全部类型" (all types) both list the
# same events, just interleaved with rows from other sheets, so the same
# seven event counts need bumping on each of the two sheets at their
# respective row numbers.

$wb = $excel.ActiveWorkbook

$exhibition = $wb.Worksheets.Item("展览")
$allTypes   = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$exhibitionUpdates = @{
    "F4"  = 303
    "F6"  = 165
    "F7"  = 309
    "F9"  = 2108
    "F11" = 5096
    "F12" = 109
    "F13" = 351
}

foreach ($addr in $exhibitionUpdates.Keys) {
    $exhibition.Range($addr).Value = $exhibitionUpdates[$addr]
}

# Sheet "全部类型": row -> new F value
$allTypesUpdates = @{
    "F5"  = 303
    "F7"  = 165
    "F8"  = 309
    "F12" = 2108
    "F14" = 5096
    "F15" = 109
    "F16" = 351
}

foreach ($addr in $allTypesUpdates.Keys) {
    $allTypes.Range($addr).Value = $allTypesUpdates[$addr]
}

$wb.Save()
